$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = [double]"0.9703253704915283"
$ws.Range("E3").Value = [double]"0.9703253704915283"

$ws.Range("D4").Value = [double]"0.9999999999957594"
$ws.Range("E4").Value = [double]"0.9999999999957594"

$ws.Range("D5").Value = [double]"1.234120095728359E-05"
$ws.Range("E5").Value = [double]"1.234120095728359E-05"

$ws.Range("D6").Value = [double]"1.459515254024199E-12"
$ws.Range("E6").Value = [double]"1.459515254024199E-12"

$ws.Range("F7").Value = [double]"11.73240280151367"
